$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.352.33'
$ws.Range('E2').Value = '  +7.71%  '
$ws.Range('D3').Value = '3.401.06'
$ws.Range('E3').Value = '  +4.85%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value2 = "'411.67"
$ws.Range('E5').Value = '  +3.69%  '
$ws.Range('D6').Value2 = "'122.15"
$ws.Range('E6').Value = '  +12.97%  '
$ws.Range('D7').Value = '3.393.72'
$ws.Range('E7').Value = '  +4.77%  '
$ws.Range('D8').Value2 = "'0.578"
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value2 = "'0.642"
$ws.Range('E10').Value = '  +3.49%  '
$ws.Range('D11').Value2 = "'0.113"
$ws.Range('E11').Value = '  +17.43%  '
$ws.Range('D12').Value2 = "'41.39"
$ws.Range('E12').Value = '  +5.24%  '
$ws.Range('E13').Value = '  -0.68%  '
$ws.Range('D14').Value = '3.943.34'
$ws.Range('E14').Value = '  +5.08%  '
$ws.Range('D15').Value2 = "'8.41"
$ws.Range('E15').Value = '  +1.03%  '
$ws.Range('D16').Value2 = "'19.54"
$ws.Range('E16').Value = '  +3.26%  '
$ws.Range('D17').Value = '3.447.81'
$ws.Range('E17').Value = '  +6.16%  '
$ws.Range('D18').Value = '61.277.23'
$ws.Range('E18').Value = '  +7.88%  '
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('D20').Value2 = "'10.85"
$ws.Range('E20').Value = '  -1.97%  '
$ws.Range('E21').Value = '  +7.05%  '
$ws.Range('D22').Value2 = "'3.34"
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').Value2 = "'12.86"
$ws.Range('E23').Value = '  -1.27%  '
$ws.Range('D24').Value2 = "'298.59"
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('D25').Value2 = "'76.15"
$ws.Range('E25').Value = '  +2.31%  '
$ws.Range('D26').Value2 = "'3.13"
$ws.Range('E26').Value = '  -1.62%  '
$ws.Range('D27').Value2 = "'30.80"
$ws.Range('E27').Value = '  +9.51%  '
$ws.Range('D28').Value2 = "'8.19"
$ws.Range('E28').Value = '  +13.19%  '
$ws.Range('D29').Value2 = "'7.67"
$ws.Range('E29').Value = '  -1.61%  '
$ws.Range('E30').Value = '  -2.12%  '
$ws.Range('E32').Value = '  +4.67%  '
$ws.Range('D33').Value2 = "'42.51"
$ws.Range('E33').Value = '  +1.09%  '
$ws.Range('D34').Value2 = "'11.42"
$ws.Range('E34').Value = '  +2.13%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').Value2 = "'2.50"
$ws.Range('E36').Value = '  +17.33%  '
$ws.Range('D37').Value2 = "'0.0481"
$ws.Range('E37').Value = '  -0.94%  '
$ws.Range('D38').Value2 = "'52.23"
$ws.Range('E38').Value = '  +1.66%  '
$ws.Range('E39').Value = '  +2.73%  '
$ws.Range('D40').Value2 = "'0.998"
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('E42').Value = '  +4.66%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').Value2 = "'0.122"
$ws.Range('E43').Value = '  +0.50%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value2 = "'134.32"
$ws.Range('E44').Value = '  -1.74%  '
$ws.Range('D45').Value2 = "'17.39"
$ws.Range('E45').Value = '  +3.38%  '
$ws.Range('D46').Value2 = "'3.94"
$ws.Range('E46').Value = '  -0.75%  '
$ws.Range('E47').Value = '  +2.05%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value2 = "'22.02"
$ws.Range('E48').Value = '  -2.36%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').Value2 = "'2.20"
$ws.Range('E49').Value = '  -3.26%  '
$ws.Range('D50').Value = '2.198.01'
$ws.Range('E50').Value = '  +2.07%  '
$ws.Range('D51').Value = '3.746.03'
$ws.Range('E51').Value = '  +5.21%  '
